$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder field text: 9/22/18 -> 9/27/18
#    (slide master + every slide layout's "Date Placeholder" shape)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "9/22/18") {
            $tr.Text = "9/27/18"
        }
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapes = $layout.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "9/22/18") {
                $tr.Text = "9/27/18"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 2 content updates
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Title: "Were you here Last Saturday?" -> "Welcome to Study Saturday!"
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Welcome to Study Saturday!"

$body = $s2.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "If not, then follow these steps:" -> "Follow these steps:"
$para1 = $body.Paragraphs(1)
$fullPara1 = $body.Characters($para1.Start, $para1.Length)
$fullPara1.Text = "Follow these steps:"

# Paragraph 2 bullet: " Clone repo: " -> " Fork & Clone " + "repo: " (two runs)
$para2 = $body.Paragraphs(2)
$oldRun = $body.Characters($para2.Start, 13)
$oldRun.Text = "repo: "
$para2.InsertBefore(" Fork & Clone ") | Out-Null
